# Intermediate step for the QUERY function: wrap the original QUERY(...) call
# with INDEX(...) to pull out individual pieces of the (still unsupported)
# QUERY result, across a 2x2 block of cells (A1:A2 on top, A3:A4 below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Formula = '=INDEX(QUERY("querydefine.xlsx","pers",48,36,"Tom B."),1,1)'
$ws.Range("A2").Formula = '=INDEX(QUERY("querydefine.xlsx","pers",48,36,"Tom B."),1,2)'
$ws.Range("A3").Formula = '=INDEX(QUERY("querydefine.xlsx","pers",48,36,"Tom B."),2,1)'
$ws.Range("A4").Formula = '=INDEX(QUERY("querydefine.xlsx","pers",48,36,"Tom B."),2,2)'

# Move the selection down to A5, like in the edited workbook.
$ws.Range("A5").Select()

# The window was scrolled down slightly in the saved view.
$excel.ActiveWindow.Top = $excel.ActiveWindow.Top
$wb.Windows.Item(1).Top = $wb.Windows.Item(1).Top
